$p = $ppt.ActivePresentation

# Duplicate the last slide ("User Account Removal") to create the new
# "Administrator Login" user-story slide at the end of the deck.
$srcSlide = $p.Slides.Item($p.Slides.Count)
$newSlide = $srcSlide.Duplicate()
$s = $p.Slides.Item($newSlide.SlideIndex)

# --- Shape 1 (id=4, "Rectangle 3") : "Story 41" -> "Story 42" -----------
$shpStory = $s.Shapes.Item(1)
$trStory = $shpStory.TextFrame.TextRange
$trStory.Text = "Story "
$trStory.InsertAfter("42") | Out-Null

# --- Shape 2 (id=6, "Rectangle 5") : title ------------------------------
$shpTitle = $s.Shapes.Item(2)
$shpTitle.TextFrame.TextRange.Text = "Administrator Login"

# --- Shape 3 (id=7, "Rectangle 6") : "As an administrator..." ----------
$shpAs = $s.Shapes.Item(3)
$trAs = $shpAs.TextFrame.TextRange
$trAs.Text = "As an administrator I want to be "
$trAs.InsertAfter("able to login to my account to make changes to the website") | Out-Null

# --- Shape 4 (id=8, "Rectangle 7") : Acceptance Criteria bullets -------
$shpAc = $s.Shapes.Item(4)
$trAc = $shpAc.TextFrame.TextRange
# Drop the 3rd bullet ("Confirmation of action dialog") by merging its
# paragraph away (deleting a non-trailing paragraph keeps the XML clean).
$trAc.Paragraphs(3, 2).Delete()
# Paragraph 2 ("Button for removing user accounts" originally) becomes
# the first bullet.
$trAc.Paragraphs(2, 1).Text = "Access to separate page from main users for logging in"
# Paragraph 3 (was "Synchronisation with user database") becomes the
# second bullet, built out of three runs to match the edited phrasing.
$p3 = $trAc.Paragraphs(3, 1)
$p3.Text = "Relevant "
$p3.InsertAfter("synchronisation") | Out-Null
$p3.InsertAfter(" with database") | Out-Null

# --- Shape 6 (id=12, "Rectangle 11") : Priority "C" -> "M" --------------
$shpPriority = $s.Shapes.Item(6)
$shpPriority.TextFrame.TextRange.Paragraphs(2, 1).Text = "M"
